# Update profit files after running on 2026-02-21
# Append a new data row (row 89) to Sheet1 with the day's profit figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 89

# Column A holds the date as literal text (matching the existing column's
# style), not an auto-converted date serial number. Briefly force a Text
# number format so the "MM/DD/YYYY"-looking string is stored verbatim, then
# clear the format again so the cell ends up with no explicit style applied
# (matching the rest of the column).
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "02/21/2026"
$ws.Range("A$row").ClearFormats()

$ws.Range("B$row").Value = 9374.48
$ws.Range("C$row").Value = 0.2438417473332335
$ws.Range("D$row").Value = 0.7561582526667665
$ws.Range("E$row").Value = -325.26
$ws.Range("F$row").Value = -34.79
$ws.Range("G$row").Value = -23903.16
$ws.Range("H$row").Value = -77.13
$ws.Range("I$row").Value = -1117.14
$ws.Range("J$row").Value = -32.83
$ws.Range("K$row").Value = -25020.12
$ws.Range("L$row").Value = -72.73999999999999
